# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
#
# The "K" column (column G) is recomputed for every data row (rows 2-28)
# and rewritten with its new value. All other columns are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2  = 0
    3  = 0
    4  = 0
    5  = 0
    6  = 1
    7  = 3
    8  = 0
    9  = 1
    10 = 0
    11 = 0
    12 = 0
    13 = 2
    14 = 1
    15 = 0
    16 = 0
    17 = 1
    18 = 0
    19 = 1
    20 = 0
    21 = 1
    22 = 0
    23 = 1
    24 = 2
    25 = 1
    26 = 0
    27 = 3
    28 = 1
}

foreach ($row in $newK.Keys) {
    $ws.Cells.Item($row, 7).Value = $newK[$row]
}
